$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "UN" (B) and date-header (C) columns.
# This shifts old column B -> D and old column C -> E, matching the target layout:
#   B = newest week (Jun_17), C = Jun_15, D = Jun_13 (old B), E = Jun_10 (old C)
$ws.Columns("B:C").Insert()

# New header values for the newly inserted columns (set C before B so the
# shared-string table gets "Jun_15" allocated before "Jun_17", matching the
# target order).
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill the newly inserted B:C columns (rows 2-27) with "UN", matching column D/E.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Re-apply the explicit "8.0 characters" custom width to columns C, D and E
# (ColumnWidth 7.1666... round-trips to a stored OOXML width of exactly 8.0).
$ws.Columns("C:E").ColumnWidth = 7.166666666666667
